{"js": "// Replace the date line and each two-digit multiplication problem in the table\n// with its updated text, matched via exact (case-sensitive) search.\nconst replacements = [\n  { find: \"2023-10-09 Monday\", replace: \"2023-10-10 Tuesday\" },\n  { find: \"58\u00d794=\", replace: \"37\u00d749=\" },\n  { find: \"21\u00d729=\", replace: \"69\u00d749=\" },\n  { find: \"12\u00d711=\", replace: \"49\u00d782=\" },\n  { find: \"40\u00d772=\", replace: \"45\u00d721=\" },\n  { find: \"65\u00d736=\", replace: \"95\u00d736=\" },\n  { find: \"11\u00d715=\", replace: \"99\u00d740=\" },\n  { find: \"18\u00d783=\", replace: \"40\u00d795=\" },\n  { find: \"98\u00d783=\", replace: \"79\u00d754=\" },\n  { find: \"72\u00d767=\", replace: \"98\u00d771=\" },\n  { find: \"39\u00d726=\", replace: \"50\u00d785=\" },\n  { find: \"95\u00d776=\", replace: \"28\u00d725=\" },\n  { find: \"34\u00d749=\", replace: \"37\u00d728=\" },\n  { find: \"45\u00d764=\", replace: \"34\u00d738=\" },\n  { find: \"35\u00d719=\", replace: \"75\u00d779=\" },\n  { find: \"13\u00d720=\", replace: \"96\u00d742=\" },\n  { find: \"22\u00d767=\", replace: \"60\u00d713=\" },\n  { find: \"44\u00d739=\", replace: \"43\u00d744=\" },\n  { find: \"88\u00d749=\", replace: \"53\u00d778=\" },\n  { find: \"12\u00d772=\", replace: \"46\u00d771=\" },\n  { find: \"59\u00d776=\", replace: \"92\u00d767=\" },\n  { find: \"36\u00d746=\", replace: \"86\u00d724=\" },\n  { find: \"53\u00d738=\", replace: \"27\u00d730=\" },\n  { find: \"32\u00d752=\", replace: \"45\u00d782=\" },\n  { find: \"23\u00d792=\", replace: \"59\u00d766=\" },\n  { find: \"36\u00d791=\", replace: \"81\u00d776=\" },\n];\n\nconst body = context.document.body;\nconst searchResults = replacements.map((pair) =>\n  body.search(pair.find, { matchCase: true, matchWholeWord: false })\n);\nsearchResults.forEach((r) => r.load(\"items\"));\n\nawait context.sync();\n\nsearchResults.forEach((results, idx) => {\n  const { find, replace } = replacements[idx];\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${find}`);\n  }\n  results.items.forEach((range) => {\n    range.insertText(replace, Word.InsertLocation.replace);\n  });\n});\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = '2023-10-09 Monday'; New = '2023-10-10 Tuesday' }\n    @{ Old = '58\u00d794='; New = '37\u00d749=' }\n    @{ Old = '21\u00d729='; New = '69\u00d749=' }\n    @{ Old = '12\u00d711='; New = '49\u00d782=' }\n    @{ Old = '40\u00d772='; New = '45\u00d721=' }\n    @{ Old = '65\u00d736='; New = '95\u00d736=' }\n    @{ Old = '11\u00d715='; New = '99\u00d740=' }\n    @{ Old = '18\u00d783='; New = '40\u00d795=' }\n    @{ Old = '98\u00d783='; New = '79\u00d754=' }\n    @{ Old = '72\u00d767='; New = '98\u00d771=' }\n    @{ Old = '39\u00d726='; New = '50\u00d785=' }\n    @{ Old = '95\u00d776='; New = '28\u00d725=' }\n    @{ Old = '34\u00d749='; New = '37\u00d728=' }\n    @{ Old = '45\u00d764='; New = '34\u00d738=' }\n    @{ Old = '35\u00d719='; New = '75\u00d779=' }\n    @{ Old = '13\u00d720='; New = '96\u00d742=' }\n    @{ Old = '22\u00d767='; New = '60\u00d713=' }\n    @{ Old = '44\u00d739='; New = '43\u00d744=' }\n    @{ Old = '88\u00d749='; New = '53\u00d778=' }\n    @{ Old = '12\u00d772='; New = '46\u00d771=' }\n    @{ Old = '59\u00d776='; New = '92\u00d767=' }\n    @{ Old = '36\u00d746='; New = '86\u00d724=' }\n    @{ Old = '53\u00d738='; New = '27\u00d730=' }\n    @{ Old = '32\u00d752='; New = '45\u00d782=' }\n    @{ Old = '23\u00d792='; New = '59\u00d766=' }\n    @{ Old = '36\u00d791='; New = '81\u00d776=' }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    # wdReplaceAll = 2, wdFindContinue = 1 (wrap within the searched range)\n    $found = $find.Execute($pair.Old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n    if (-not $found) {\n        throw \"Could not find text to replace: $($pair.Old)\"\n    }\n}"}
